$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cell H1 - style it like the other header cells (copy style from G1)
$ws.Range("H1").Value = "Save"
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122) # xlPasteFormats

# Save flag values for rows 2-19
$saveValues = @{
    2  = 0
    3  = 0
    4  = 1
    5  = 0
    6  = 0
    7  = 0
    8  = 1
    9  = 0
    10 = 0
    11 = 0
    12 = 0
    13 = 0
    14 = 1
    15 = 1
    16 = 1
    17 = 0
    18 = 1
    19 = 1
}

foreach ($row in $saveValues.Keys) {
    $ws.Cells.Item($row, 8).Value = $saveValues[$row]
}
